# Auto-generated edit script
# Commit: Add data for 2025-12-31
# Applies 194 individual cell updates (2025 YTD crime counts) across 27 sheets

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Cells.Item(2, 7).Value = 93  # G2
$ws.Cells.Item(3, 4).Value = 140  # D3
$ws.Cells.Item(3, 9).Value = 200  # I3
$ws.Cells.Item(3, 10).Value = 243  # J3
$ws.Cells.Item(3, 12).Value = 262  # L3
$ws.Cells.Item(9, 3).Value = 509  # C9
$ws.Cells.Item(9, 4).Value = 446  # D9
$ws.Cells.Item(9, 5).Value = 516  # E9
$ws.Cells.Item(9, 6).Value = 589  # F9
$ws.Cells.Item(9, 7).Value = 450  # G9
$ws.Cells.Item(9, 8).Value = 483  # H9
$ws.Cells.Item(9, 9).Value = 520  # I9
$ws.Cells.Item(9, 10).Value = 439  # J9
$ws.Cells.Item(9, 12).Value = 463  # L9
$ws.Cells.Item(10, 2).Value = 1426  # B10
$ws.Cells.Item(10, 3).Value = 1673  # C10
$ws.Cells.Item(10, 4).Value = 1894  # D10
$ws.Cells.Item(10, 5).Value = 2338  # E10
$ws.Cells.Item(10, 6).Value = 2225  # F10
$ws.Cells.Item(10, 7).Value = 928  # G10
$ws.Cells.Item(10, 8).Value = 640  # H10
$ws.Cells.Item(10, 10).Value = 766  # J10
$ws.Cells.Item(10, 11).Value = 711  # K10
$ws.Cells.Item(10, 12).Value = 702  # L10
$ws.Cells.Item(11, 2).Value = 1965  # B11
$ws.Cells.Item(11, 3).Value = 2348  # C11
$ws.Cells.Item(11, 4).Value = 2588  # D11
$ws.Cells.Item(11, 5).Value = 3095  # E11
$ws.Cells.Item(11, 6).Value = 3072  # F11
$ws.Cells.Item(11, 7).Value = 1630  # G11
$ws.Cells.Item(11, 8).Value = 1418  # H11
$ws.Cells.Item(11, 9).Value = 1757  # I11
$ws.Cells.Item(11, 10).Value = 1610  # J11
$ws.Cells.Item(11, 11).Value = 1672  # K11
$ws.Cells.Item(11, 12).Value = 1594  # L11

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Cells.Item(3, 9).Value = 9  # I3
$ws.Cells.Item(7, 4).Value = 34  # D7
$ws.Cells.Item(7, 6).Value = 54  # F7
$ws.Cells.Item(7, 10).Value = 34  # J7
$ws.Cells.Item(8, 4).Value = 49  # D8
$ws.Cells.Item(8, 6).Value = 138  # F8
$ws.Cells.Item(9, 4).Value = 100  # D9
$ws.Cells.Item(9, 6).Value = 204  # F9
$ws.Cells.Item(9, 9).Value = 98  # I9
$ws.Cells.Item(9, 10).Value = 85  # J9

$ws = $wb.Worksheets.Item('Chatham')
$ws.Cells.Item(9, 6).Value = 34  # F9
$ws.Cells.Item(9, 7).Value = 15  # G9
$ws.Cells.Item(9, 10).Value = 20  # J9
$ws.Cells.Item(10, 6).Value = 60  # F10
$ws.Cells.Item(10, 7).Value = 40  # G10
$ws.Cells.Item(10, 10).Value = 42  # J10

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Cells.Item(9, 12).Value = 36  # L9
$ws.Cells.Item(10, 12).Value = 99  # L10

$ws = $wb.Worksheets.Item('Loop')
$ws.Cells.Item(8, 10).Value = 64  # J8
$ws.Cells.Item(9, 5).Value = 703  # E9
$ws.Cells.Item(9, 6).Value = 565  # F9
$ws.Cells.Item(9, 7).Value = 168  # G9
$ws.Cells.Item(9, 12).Value = 96  # L9
$ws.Cells.Item(10, 5).Value = 794  # E10
$ws.Cells.Item(10, 6).Value = 656  # F10
$ws.Cells.Item(10, 7).Value = 257  # G10
$ws.Cells.Item(10, 10).Value = 255  # J10
$ws.Cells.Item(10, 12).Value = 197  # L10

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Cells.Item(6, 5).Value = 11  # E6
$ws.Cells.Item(7, 5).Value = 23  # E7

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Cells.Item(6, 2).Value = 17  # B6
$ws.Cells.Item(7, 2).Value = 27  # B7

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Cells.Item(6, 5).Value = 19  # E6
$ws.Cells.Item(6, 8).Value = 10  # H6
$ws.Cells.Item(7, 6).Value = 58  # F7
$ws.Cells.Item(8, 5).Value = 66  # E8
$ws.Cells.Item(8, 6).Value = 100  # F8
$ws.Cells.Item(8, 8).Value = 26  # H8

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Cells.Item(5, 5).Value = 23  # E5
$ws.Cells.Item(8, 5).Value = 130  # E8
$ws.Cells.Item(8, 6).Value = 167  # F8
$ws.Cells.Item(8, 12).Value = 62  # L8
$ws.Cells.Item(19, 6).Value = 60  # F19
$ws.Cells.Item(19, 7).Value = 40  # G19
$ws.Cells.Item(19, 10).Value = 42  # J19
$ws.Cells.Item(28, 2).Value = 111  # B28
$ws.Cells.Item(28, 5).Value = 98  # E28
$ws.Cells.Item(28, 6).Value = 138  # F28
$ws.Cells.Item(28, 10).Value = 69  # J28
$ws.Cells.Item(32, 4).Value = 100  # D32
$ws.Cells.Item(32, 6).Value = 204  # F32
$ws.Cells.Item(32, 9).Value = 98  # I32
$ws.Cells.Item(32, 10).Value = 85  # J32
$ws.Cells.Item(36, 12).Value = 99  # L36
$ws.Cells.Item(43, 4).Value = 21  # D43
$ws.Cells.Item(43, 6).Value = 18  # F43
$ws.Cells.Item(45, 8).Value = 20  # H45
$ws.Cells.Item(48, 2).Value = 15  # B48
$ws.Cells.Item(49, 2).Value = 7  # B49
$ws.Cells.Item(49, 3).Value = 16  # C49
$ws.Cells.Item(50, 2).Value = 27  # B50
$ws.Cells.Item(51, 3).Value = 16  # C51
$ws.Cells.Item(52, 7).Value = 26  # G52
$ws.Cells.Item(52, 8).Value = 18  # H52
$ws.Cells.Item(53, 5).Value = 794  # E53
$ws.Cells.Item(53, 6).Value = 656  # F53
$ws.Cells.Item(53, 7).Value = 257  # G53
$ws.Cells.Item(53, 10).Value = 255  # J53
$ws.Cells.Item(53, 12).Value = 197  # L53
$ws.Cells.Item(54, 5).Value = 20  # E54
$ws.Cells.Item(61, 2).Value = 22  # B61
$ws.Cells.Item(61, 12).Value = 2  # L61
$ws.Cells.Item(64, 9).Value = 7  # I64
$ws.Cells.Item(65, 5).Value = 66  # E65
$ws.Cells.Item(65, 6).Value = 100  # F65
$ws.Cells.Item(65, 8).Value = 26  # H65
$ws.Cells.Item(74, 4).Value = 97  # D74
$ws.Cells.Item(74, 5).Value = 82  # E74
$ws.Cells.Item(74, 7).Value = 37  # G74
$ws.Cells.Item(76, 7).Value = 53  # G76
$ws.Cells.Item(76, 12).Value = 60  # L76
$ws.Cells.Item(77, 12).Value = 62  # L77
$ws.Cells.Item(80, 11).Value = 17  # K80
$ws.Cells.Item(84, 12).Value = 5  # L84
$ws.Cells.Item(85, 9).Value = 4  # I85
$ws.Cells.Item(85, 12).Value = 5  # L85
$ws.Cells.Item(86, 8).Value = 9  # H86
$ws.Cells.Item(89, 2).Value = 28  # B89
$ws.Cells.Item(99, 2).Value = 1965  # B99
$ws.Cells.Item(99, 3).Value = 2348  # C99
$ws.Cells.Item(99, 4).Value = 2588  # D99
$ws.Cells.Item(99, 5).Value = 3095  # E99
$ws.Cells.Item(99, 6).Value = 3072  # F99
$ws.Cells.Item(99, 7).Value = 1630  # G99
$ws.Cells.Item(99, 8).Value = 1418  # H99
$ws.Cells.Item(99, 9).Value = 1757  # I99
$ws.Cells.Item(99, 10).Value = 1610  # J99
$ws.Cells.Item(99, 11).Value = 1672  # K99
$ws.Cells.Item(99, 12).Value = 1594  # L99

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Cells.Item(7, 2).Value = 21  # B7
$ws.Cells.Item(8, 2).Value = 28  # B8

$ws = $wb.Worksheets.Item('Sheffield & DePaul')
$ws.Cells.Item(6, 11).Value = 9  # K6
$ws.Cells.Item(7, 11).Value = 17  # K7

$ws = $wb.Worksheets.Item('Englewood')
$ws.Cells.Item(3, 10).Value = 12  # J3
$ws.Cells.Item(8, 2).Value = 71  # B8
$ws.Cells.Item(8, 5).Value = 65  # E8
$ws.Cells.Item(8, 6).Value = 75  # F8
$ws.Cells.Item(9, 2).Value = 111  # B9
$ws.Cells.Item(9, 5).Value = 98  # E9
$ws.Cells.Item(9, 6).Value = 138  # F9
$ws.Cells.Item(9, 10).Value = 69  # J9

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Cells.Item(6, 8).Value = 11  # H6
$ws.Cells.Item(7, 8).Value = 20  # H7

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Cells.Item(2, 7).Value = 6  # G2
$ws.Cells.Item(3, 12).Value = 6  # L3
$ws.Cells.Item(10, 7).Value = 53  # G10
$ws.Cells.Item(10, 12).Value = 60  # L10

$ws = $wb.Worksheets.Item('River North')
$ws.Cells.Item(3, 4).Value = 4  # D3
$ws.Cells.Item(5, 7).Value = 13  # G5
$ws.Cells.Item(6, 5).Value = 75  # E6
$ws.Cells.Item(7, 4).Value = 97  # D7
$ws.Cells.Item(7, 5).Value = 82  # E7
$ws.Cells.Item(7, 7).Value = 37  # G7

$ws = $wb.Worksheets.Item('United Center')
$ws.Cells.Item(6, 8).Value = 3  # H6
$ws.Cells.Item(8, 8).Value = 9  # H8

$ws = $wb.Worksheets.Item('Little Village')
$ws.Cells.Item(5, 3).Value = 4  # C5
$ws.Cells.Item(7, 3).Value = 16  # C7

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Cells.Item(7, 7).Value = 17  # G7
$ws.Cells.Item(7, 8).Value = 8  # H7
$ws.Cells.Item(8, 7).Value = 26  # G8
$ws.Cells.Item(8, 8).Value = 18  # H8

$ws = $wb.Worksheets.Item('North Center')
$ws.Cells.Item(5, 9).Value = 1  # I5
$ws.Cells.Item(7, 9).Value = 7  # I7

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Cells.Item(5, 5).Value = 13  # E5
$ws.Cells.Item(6, 5).Value = 20  # E6

$ws = $wb.Worksheets.Item('Roseland')
$ws.Cells.Item(9, 12).Value = 25  # L9
$ws.Cells.Item(10, 12).Value = 62  # L10

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Cells.Item(6, 2).Value = 6  # B6
$ws.Cells.Item(6, 3).Value = 14  # C6
$ws.Cells.Item(7, 2).Value = 7  # B7
$ws.Cells.Item(7, 3).Value = 16  # C7

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Cells.Item(5, 12).Value = 1  # L5
$ws.Cells.Item(7, 12).Value = 5  # L7

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Cells.Item(6, 2).Value = 13  # B6
$ws.Cells.Item(7, 2).Value = 15  # B7

$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Cells.Item(3, 12).Value = 2  # L3
$ws.Cells.Item(4, 9).Value = 2  # I4
$ws.Cells.Item(6, 9).Value = 4  # I6
$ws.Cells.Item(6, 12).Value = 5  # L6

$ws = $wb.Worksheets.Item('Austin')
$ws.Cells.Item(9, 5).Value = 73  # E9
$ws.Cells.Item(9, 6).Value = 114  # F9
$ws.Cells.Item(9, 12).Value = 27  # L9
$ws.Cells.Item(10, 5).Value = 130  # E10
$ws.Cells.Item(10, 6).Value = 167  # F10
$ws.Cells.Item(10, 12).Value = 62  # L10

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Cells.Item(5, 4).Value = 5  # D5
$ws.Cells.Item(6, 6).Value = 15  # F6
$ws.Cells.Item(7, 4).Value = 21  # D7
$ws.Cells.Item(7, 6).Value = 18  # F7
